# Refresh the cryptos price list (GitHub Actions style scheduled update).
# All Price/Volume columns in this sheet are stored as plain text
# (t="inlineStr" in the original OOXML), even when a value happens to look
# like a decimal number (e.g. "65.12"). Excel's COM Range.Value setter
# auto-detects such strings as numbers, so for any replacement price that
# parses as a plain number we briefly force the cell to Text format ("@")
# before assigning it, then ClearFormats() right after so the cell ends up
# with no explicit style - matching the unstyled text cells in the source
# file - while keeping the stored value as text, not a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.917.74"
$ws.Range("E2").Value = "  +3.89%  "
$ws.Range("D3").Value = "2.225.03"
$ws.Range("E3").Value = "  +3.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.20"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("E6").Value = "  +1.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.12"
$ws.Range("E7").Value = "  +1.59%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +2.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0886"
$ws.Range("E10").Value = "  +3.00%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "2.555.89"
$ws.Range("E12").Value = "  +3.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "16.13"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.43"
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.825"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.67"
$ws.Range("E16").Value = "  +2.20%  "
$ws.Range("D17").Value = "2.224.50"
$ws.Range("E17").Value = "  +2.97%  "
$ws.Range("D18").Value = "40.813.45"
$ws.Range("E18").Value = "  +3.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.07"
$ws.Range("E19").Value = "  +3.17%  "
$ws.Range("D20").Value = "0.0₃0912"
$ws.Range("E20").Value = "  +6.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.17"
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "254.65"
$ws.Range("E22").Value = "  +10.00%  "
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("E25").Value = "  -8.00%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.77"
$ws.Range("E26").Value = "  +2.49%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "173.25"
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.145"
$ws.Range("E28").Value = "  +4.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.39"
$ws.Range("E29").Value = "  +2.40%  "
$ws.Range("E30").Value = "  +2.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.84"
$ws.Range("E31").Value = "  +5.52%  "
$ws.Range("E32").Value = "  +1.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.69"
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("E34").Value = "  +1.97%  "
$ws.Range("E35").Value = "  +2.34%  "
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.83"
$ws.Range("E37").Value = "  +6.21%  "
$ws.Range("E38").Value = "  +3.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.91"
$ws.Range("E40").Value = "  +16.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0237"
$ws.Range("E41").Value = "  +2.97%  "
$ws.Range("B42").Value = "TerraClassic"
$ws.Range("C42").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.000232"
$ws.Range("E42").Value = "  +53.96%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.77"
$ws.Range("E43").Value = "  +10.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.35"
$ws.Range("E44").Value = "  -1.47%  "
$ws.Range("E45").Value = "  +5.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.67"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").Value = "1.517.07"
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0947"
$ws.Range("E48").Value = "  +2.44%  "
$ws.Range("E49").Value = "  +2.11%  "
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.13"
$ws.Range("E51").Value = "  +12.65%  "

# Restore default (unstyled) formatting for cells forced to text above
$ws.Range("D4").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D51").ClearFormats()
